# Natmi following Dr Hou advice
# Update the LR-pairs table: recompute existing rows (ECs->ECs self row and
# ECs->FAPs row) and add two new target-cluster rows (M1 and sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Ligand=Fgf9, Receptor=Fgfr2, Target=ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf9"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.883656666666667
$ws.Range("H2").Value = 5.650970000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.493831
$ws.Range("N2").Value = 1.481493
$ws.Range("O2").Value = 0.1121895146192186
$ws.Range("P2").Value = 0.1134277760249069
$ws.Range("Q2").Value = 0.9302080553566667
$ws.Range("R2").Value = 8.371872498210001
$ws.Range("S2").Value = 0.1121895146192186
$ws.Range("T2").Value = 0.1134277760249069

# Row 3: Sending=ECs, Ligand=Fgf9, Receptor=Fgfr2, Target=FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf9"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.883656666666667
$ws.Range("H3").Value = 5.650970000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.763360333333333
$ws.Range("N3").Value = 11.290081
$ws.Range("O3").Value = 0.85496773012202
$ws.Range("P3").Value = 0.8644042050627692
$ws.Range("Q3").Value = 7.088878780952223
$ws.Range("R3").Value = 63.79990902857001
$ws.Range("S3").Value = 0.85496773012202
$ws.Range("T3").Value = 0.8644042050627692

# Row 4 (new): Sending=ECs, Ligand=Fgf9, Receptor=Fgfr2, Target=M1
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf9"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.883656666666667
$ws.Range("H4").Value = 5.650970000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.0004073333333333333
$ws.Range("N4").Value = 0.001222
$ws.Range("O4").Value = 0.00009253880164447964
$ws.Range("P4").Value = 0.00009356017362379453
$ws.Range("Q4").Value = 0.000767276148888889
$ws.Range("R4").Value = 0.006905485340000002
$ws.Range("S4").Value = 0.00009253880164447964
$ws.Range("T4").Value = 0.00009356017362379453

# Row 5 (new): Sending=ECs, Ligand=Fgf9, Receptor=Fgfr2, Target=sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf9"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.883656666666667
$ws.Range("H5").Value = 5.650970000000001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1441585
$ws.Range("N5").Value = 0.288317
$ws.Range("O5").Value = 0.03275021645711715
$ws.Range("P5").Value = 0.02207445873870014
$ws.Range("Q5").Value = 0.2715451195816667
$ws.Range("R5").Value = 1.62927071749
$ws.Range("S5").Value = 0.03275021645711715
$ws.Range("T5").Value = 0.02207445873870014
